# Update the worksheet date and the 25 two-digit multiplication problems
# to the "output generated at 503736d" values.

$d = $word.ActiveDocument

# Mapping of old text -> new text (each old value is unique in the document).
$replacements = @(
    @{ Old = "2025-01-19 Sunday"; New = "2025-01-20 Monday" },
    @{ Old = "99×73="; New = "90×27=" },
    @{ Old = "78×84="; New = "69×58=" },
    @{ Old = "22×95="; New = "91×29=" },
    @{ Old = "90×96="; New = "93×23=" },
    @{ Old = "74×35="; New = "39×83=" },
    @{ Old = "56×37="; New = "98×26=" },
    @{ Old = "66×62="; New = "12×70=" },
    @{ Old = "49×47="; New = "11×45=" },
    @{ Old = "53×17="; New = "76×30=" },
    @{ Old = "56×96="; New = "19×26=" },
    @{ Old = "43×68="; New = "43×54=" },
    @{ Old = "66×76="; New = "39×66=" },
    @{ Old = "14×22="; New = "64×69=" },
    @{ Old = "64×61="; New = "86×38=" },
    @{ Old = "75×33="; New = "31×27=" },
    @{ Old = "54×23="; New = "12×20=" },
    @{ Old = "68×97="; New = "50×90=" },
    @{ Old = "76×43="; New = "78×94=" },
    @{ Old = "84×52="; New = "47×62=" },
    @{ Old = "80×45="; New = "58×14=" },
    @{ Old = "93×31="; New = "61×95=" },
    @{ Old = "97×87="; New = "70×22=" },
    @{ Old = "44×16="; New = "20×35=" },
    @{ Old = "36×38="; New = "89×53=" },
    @{ Old = "91×32="; New = "99×89=" }
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute(
        $r.Old,    # FindText
        $true,     # MatchCase
        $false,    # MatchWholeWord
        $false,    # MatchWildcards
        $false,    # MatchSoundsLike
        $false,    # MatchAllWordForms
        $true,     # Forward
        1,         # Wrap (wdFindContinue)
        $false,    # Format
        $r.New,    # ReplaceWith
        2          # Replace (wdReplaceAll)
    )
    if (-not $found) {
        Write-Host "WARNING: text not found for replacement:" $r.Old
    }
}
